$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 500
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1500
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 1567
$ws.Range("N111").ClearContents()

# Hunk 1: ARM!row30
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 4966.3335
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 5949.5
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 5949.5
$ws.Range("M30").Value = -2850
$ws.Range("N30").Value = -6249.5

# Hunk 2: ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 64044.125
$ws.Range("I45").Value = 72836.21000000001
$ws.Range("K45").Value = 72836.21000000001
$ws.Range("M45").Value = -72459.21000000001

# Hunk 3: ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1924.2858
$ws.Range("I97").Value = 1761.6666
$ws.Range("K97").Value = 1761.6666
$ws.Range("M97").Value = -1265.6666

# Hunk 4: ARM!row119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 21333
$ws.Range("J119").Value = 21333
$ws.Range("L119").Value = 21333
$ws.Range("N119").Value = -31009

# Hunk 5: ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1141.7037
$ws.Range("I122").Value = 1154.8462
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 3464.5386
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -1014.5386
$ws.Range("N122").Value = -7300

# Hunk 6: ARM!row135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 54271.5
$ws.Range("J135").Value = 54271.5
$ws.Range("L135").Value = 54271.5
$ws.Range("N135").Value = -64411.5

# Hunk 7: BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1168
$ws.Range("I86").Value = 1001
$ws.Range("J86").Value = 1376.75
$ws.Range("K86").Value = 1001
$ws.Range("L86").Value = 1376.75
$ws.Range("M86").Value = 122
$ws.Range("N86").Value = -3622.75

# Hunk 8: BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1168
$ws.Range("I89").Value = 1001
$ws.Range("J89").Value = 1376.75
$ws.Range("K89").Value = 5005
$ws.Range("L89").Value = 6883.75
$ws.Range("M89").Value = 611
$ws.Range("N89").Value = -18115.75

# Hunk 9: BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1004
$ws.Range("I99").Value = 981.8461
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 981.8461
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 516.1539
$ws.Range("N99").Value = -4096

# Hunk 10: BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 34484424
$ws.Range("I105").Value = 1600
$ws.Range("J105").Value = 76924820
$ws.Range("K105").Value = 1600
$ws.Range("L105").Value = 76924820
$ws.Range("M105").Value = 147
$ws.Range("N105").Value = -76928314

# Hunk 11: CRP!row19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2020.1
$ws.Range("I19").Value = 1337.625
$ws.Range("J19").Value = 4750
$ws.Range("K19").Value = 1337.625
$ws.Range("L19").Value = 4750
$ws.Range("M19").Value = -1167.625
$ws.Range("N19").Value = -5090

# Hunk 12: CRP!row24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 2020.1
$ws.Range("I24").Value = 1337.625
$ws.Range("J24").Value = 4750
$ws.Range("K24").Value = 1337.625
$ws.Range("L24").Value = 4750
$ws.Range("M24").Value = -1167.625
$ws.Range("N24").Value = -5090

# Hunk 13: CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5953332
$ws.Range("I31").Value = 4630423
$ws.Range("J31").Value = 8334568.5
$ws.Range("K31").Value = 4630423
$ws.Range("L31").Value = 8334568.5
$ws.Range("M31").Value = -4630128
$ws.Range("N31").Value = -8335158.5

# Hunk 14: CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5953332
$ws.Range("I34").Value = 4630423
$ws.Range("J34").Value = 8334568.5
$ws.Range("K34").Value = 4630423
$ws.Range("L34").Value = 8334568.5
$ws.Range("M34").Value = -4630221
$ws.Range("N34").Value = -8334972.5

# Hunk 15: CRP!row53
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 15500
$ws.Range("J53").Value = 15500
$ws.Range("L53").Value = 15500
$ws.Range("N53").Value = -16714

# Hunk 16: CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 62501740
$ws.Range("I99").Value = 100001620
$ws.Range("J99").Value = 1950
$ws.Range("K99").Value = 100001620
$ws.Range("L99").Value = 1950
$ws.Range("M99").Value = -100000122
$ws.Range("N99").Value = -4946

# Hunk 17: CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 62501740
$ws.Range("I126").Value = 100001620
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 300004860
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -300002390
$ws.Range("N126").Value = -10790

# Hunk 18: CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1515.3846
$ws.Range("I134").Value = 1562.5
$ws.Range("J134").Value = 1440
$ws.Range("K134").Value = 4687.5
$ws.Range("L134").Value = 4320
$ws.Range("M134").Value = -2152.5
$ws.Range("N134").Value = -9390

# Hunk 19: CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2459.3167
$ws.Range("I68").Value = 632.125
$ws.Range("J68").Value = 3123.75
$ws.Range("K68").Value = 1896.375
$ws.Range("L68").Value = 9371.25
$ws.Range("M68").Value = -1085.375
$ws.Range("N68").Value = -10993.25

# Hunk 20: CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2459.3167
$ws.Range("I71").Value = 632.125
$ws.Range("J71").Value = 3123.75
$ws.Range("K71").Value = 5689.125
$ws.Range("L71").Value = 28113.75
$ws.Range("M71").Value = -1633.125
$ws.Range("N71").Value = -36225.75

# Hunk 21: GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2161.3044
$ws.Range("I80").Value = 2307.1428
$ws.Range("J80").Value = 1934.4445
$ws.Range("K80").Value = 2307.1428
$ws.Range("L80").Value = 1934.4445
$ws.Range("M80").Value = -1309.1428
$ws.Range("N80").Value = -3930.4445

# Hunk 22: GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2161.3044
$ws.Range("I83").Value = 2307.1428
$ws.Range("J83").Value = 1934.4445
$ws.Range("K83").Value = 11535.714
$ws.Range("L83").Value = 9672.2225
$ws.Range("M83").Value = -6543.714
$ws.Range("N83").Value = -19656.2225

# Hunk 23: LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1808808.1
$ws.Range("I22").Value = 2110109.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2110109.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -2109814.5
$ws.Range("N22").Value = -1590

# Hunk 24: LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1808808.1
$ws.Range("I27").Value = 2110109.5
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 2110109.5
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -2110002.5
$ws.Range("N27").Value = -1214

# Hunk 25: LTW!row87
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 10000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 10000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -8877
$ws.Range("N87").ClearContents()

# Hunk 26: LTW!row90
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 10000
$ws.Range("I90").Value = 10000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 30000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -24384
$ws.Range("N90").ClearContents()

# Hunk 27: LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10734.412
$ws.Range("I122").Value = 12632
$ws.Range("J122").Value = 6180.2
$ws.Range("K122").Value = 37896
$ws.Range("L122").Value = 18540.6
$ws.Range("M122").Value = -35446
$ws.Range("N122").Value = -23440.6

# Hunk 28: LTW!row133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 35775.332
$ws.Range("J133").Value = 35775.332
$ws.Range("L133").Value = 35775.332
$ws.Range("N133").Value = -40835.332

# Hunk 29: WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 66774.875
$ws.Range("I122").Value = 87366.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 262099.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -259649.5
$ws.Range("N122").Value = -19900
